$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $escaped = $value -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $ws.Range('D2') '43.718.64'
Set-TextValue $ws.Range('E2') '  +5.02%  '
Set-TextValue $ws.Range('D3') '2.248.95'
Set-TextValue $ws.Range('E3') '  +2.30%  '
Set-TextValue $ws.Range('E4') '  +0.02%  '
Set-TextValue $ws.Range('D5') '229.52'
Set-TextValue $ws.Range('E5') '  +0.12%  '
Set-TextValue $ws.Range('D6') '0.629'
Set-TextValue $ws.Range('D7') '63.52'
Set-TextValue $ws.Range('E7') '  +5.62%  '
Set-TextValue $ws.Range('E8') '  +0.01%  '
Set-TextValue $ws.Range('D9') '0.432'
Set-TextValue $ws.Range('E9') '  +8.31%  '
Set-TextValue $ws.Range('E10') '  +13.19%  '
Set-TextValue $ws.Range('D11') '56.30'
Set-TextValue $ws.Range('E11') '  -1.18%  '
Set-TextValue $ws.Range('E12') '  +2.86%  '
Set-TextValue $ws.Range('D13') '25.64'
Set-TextValue $ws.Range('E13') '  +16.76%  '
Set-TextValue $ws.Range('D14') '2.579.70'
Set-TextValue $ws.Range('E14') '  +2.14%  '
Set-TextValue $ws.Range('D15') '15.58'
Set-TextValue $ws.Range('E15') '  +1.85%  '
Set-TextValue $ws.Range('D16') '5.92'
Set-TextValue $ws.Range('E16') '  +6.81%  '
Set-TextValue $ws.Range('D17') '0.819'
Set-TextValue $ws.Range('E17') '  +3.50%  '
Set-TextValue $ws.Range('D18') '2.270.87'
Set-TextValue $ws.Range('E18') '  +3.05%  '
Set-TextValue $ws.Range('D19') '43.568.68'
Set-TextValue $ws.Range('E19') '  +4.79%  '
Set-TextValue $ws.Range('D20') '0.0000101'
Set-TextValue $ws.Range('E20') '  +11.97%  '
Set-TextValue $ws.Range('D21') '73.35'
Set-TextValue $ws.Range('E21') '  +2.27%  '
Set-TextValue $ws.Range('D22') '6.00'
Set-TextValue $ws.Range('E22') '  +0.00%  '
Set-TextValue $ws.Range('D23') '253.05'
Set-TextValue $ws.Range('E23') '  +5.16%  '
Set-TextValue $ws.Range('E24') '  +0.15%  '
Set-TextValue $ws.Range('D25') '2.43'
Set-TextValue $ws.Range('E25') '  +4.08%  '
Set-TextValue $ws.Range('E26') '  +2.07%  '
Set-TextValue $ws.Range('D27') '9.87'
Set-TextValue $ws.Range('E27') '  +3.38%  '
Set-TextValue $ws.Range('D28') '171.18'
Set-TextValue $ws.Range('E28') '  +2.01%  '
Set-TextValue $ws.Range('B29') 'EthereumClassic'
Set-TextValue $ws.Range('C29') 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D29') '20.68'
Set-TextValue $ws.Range('E29') '  +5.24%  '
Set-TextValue $ws.Range('B30') 'Kaspa'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D30') '0.136'
Set-TextValue $ws.Range('E30') '  -1.40%  '
Set-TextValue $ws.Range('E31') '  +9.45%  '
Set-TextValue $ws.Range('E32') '  -3.61%  '
Set-TextValue $ws.Range('E33') '  +2.60%  '
Set-TextValue $ws.Range('D34') '0.0674'
Set-TextValue $ws.Range('E34') '  +4.88%  '
Set-TextValue $ws.Range('D35') '4.67'
Set-TextValue $ws.Range('E35') '  +2.27%  '
Set-TextValue $ws.Range('E36') '  -0.38%  '
Set-TextValue $ws.Range('D37') '3.86'
Set-TextValue $ws.Range('E37') '  +9.32%  '
Set-TextValue $ws.Range('D38') '6.65'
Set-TextValue $ws.Range('E38') '  +6.13%  '
Set-TextValue $ws.Range('D39') '2.30'
Set-TextValue $ws.Range('E39') '  -0.84%  '
Set-TextValue $ws.Range('D40') '0.0249'
Set-TextValue $ws.Range('E40') '  +4.53%  '
Set-TextValue $ws.Range('D41') '0.999'
Set-TextValue $ws.Range('E41') '  -0.20%  '
Set-TextValue $ws.Range('B42') 'InjectiveProtocol'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D42') '17.39'
Set-TextValue $ws.Range('E42') '  +8.75%  '
Set-TextValue $ws.Range('B43') 'FraxShare'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D43') '8.19'
Set-TextValue $ws.Range('E43') '  -4.57%  '
Set-TextValue $ws.Range('D44') '0.0956'
Set-TextValue $ws.Range('E44') '  +0.31%  '
Set-TextValue $ws.Range('B45') 'Aave'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D45') '96.74'
Set-TextValue $ws.Range('E45') '  +0.54%  '
Set-TextValue $ws.Range('B46') 'TerraClassic'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
Set-TextValue $ws.Range('D46') '0.000211'
Set-TextValue $ws.Range('E46') '  -13.69%  '
Set-TextValue $ws.Range('B47') 'TrustWalletToken'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D47') '1.18'
Set-TextValue $ws.Range('E47') '  -1.26%  '
Set-TextValue $ws.Range('D48') '4.34'
Set-TextValue $ws.Range('E48') '  +0.05%  '
Set-TextValue $ws.Range('B49') 'Maker'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D49') '1.449.14'
Set-TextValue $ws.Range('E49') '  +0.13%  '
Set-TextValue $ws.Range('B50') 'NEARProtocol'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D50') '2.30'
Set-TextValue $ws.Range('E50') '  +4.63%  '
Set-TextValue $ws.Range('B51') 'HuobiToken'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D51') '2.72'

$excel.CutCopyMode = 0